# Update "想去人数" (number of people interested) values (column F) across the
# four worksheets of the workbook, reflecting refreshed stats at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1129
$ws1.Range("F4").Value  = 20435
$ws1.Range("F6").Value  = 2528
$ws1.Range("F7").Value  = 780
$ws1.Range("F9").Value  = 480
$ws1.Range("F10").Value = 730
$ws1.Range("F11").Value = 269
$ws1.Range("F14").Value = 393
$ws1.Range("F15").Value = 97
$ws1.Range("F21").Value = 112

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 310
$ws2.Range("F6").Value  = 139
$ws2.Range("F9").Value  = 12
$ws2.Range("F14").Value = 120

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6073
$ws3.Range("F3").Value = 677
$ws3.Range("F4").Value = 647
$ws3.Range("F5").Value = 1374
$ws3.Range("F6").Value = 34

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6073
$ws4.Range("F3").Value  = 677
$ws4.Range("F4").Value  = 647
$ws4.Range("F5").Value  = 1374
$ws4.Range("F7").Value  = 1129
$ws4.Range("F8").Value  = 20435
$ws4.Range("F13").Value = 310
$ws4.Range("F14").Value = 2528
$ws4.Range("F15").Value = 780
$ws4.Range("F16").Value = 139
$ws4.Range("F17").Value = 34
$ws4.Range("F19").Value = 480
$ws4.Range("F20").Value = 730
$ws4.Range("F21").Value = 269
$ws4.Range("F26").Value = 12
$ws4.Range("F27").Value = 393
$ws4.Range("F28").Value = 97
$ws4.Range("F36").Value = 120
$ws4.Range("F37").Value = 120
$ws4.Range("F48").Value = 112
